$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "302.84"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.23%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "32.17"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.24%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.998"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.83%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07870"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.99%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.079"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-5.50%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.837"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.19%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.845"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.52%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9265"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.03%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1764"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.87%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07790"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "6.13%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08587"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-5.78%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03165"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "4.59%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1006"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.39%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001512"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.52%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005686"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-4.69%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2,108.96%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.462"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.52%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-6.49%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.15%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.56%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.289"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.77%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "16.87%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.84%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-1.42%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004455"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.12%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001250"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "4.21%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01734"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-0.65%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04801"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "4.49%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007476"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "7.75%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1365"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.58%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002360"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "7.82%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01043"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "9.27%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006131"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-2.57%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.04%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.003100"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-61.13%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.8205"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "9.83%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.04%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002000"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.04%"
